$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.245.68"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.907.01"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5369"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3824"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9057"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08204"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.365"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008655"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "27.270.44"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.050"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "1.073.40"
$ws.Range("E21").Value = "  -43.31%  "
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.744"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.827"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.720"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09221"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8300"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05088"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.219"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.005"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.332"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.668"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5897"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02003"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.079"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5061"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1528"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.643"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06161"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
